$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 (Ost) and R2 (Nord) to whole numbers
$ws.Range("Q2").Value = 565783
$ws.Range("R2").Value = 6956702

# Remove the time cells Z2 (Starttid) and AB2 (Sluttid) entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
